$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EA_incidence")

# Update the base rate in C6 - this cascades via formulas to C2:C11, G2:G11, K2:K11
$ws.Range("C6").Value = 225

# Update the active cell selection to match the saved view state
$ws.Range("C6").Select()
